$d = $word.ActiveDocument

# -----------------------------------------------------------------------
# Locate the "OBS: ..." paragraph, and the first blank paragraph that
# immediately follows it, without relying on hard-coded indices.
# -----------------------------------------------------------------------
$obsIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text.StartsWith("OBS:")) {
        $obsIndex = $i
        break
    }
}
if ($obsIndex -eq -1) {
    throw "Could not locate the 'OBS:' paragraph"
}

$blankIndex = -1
for ($j = $obsIndex + 1; $j -le $d.Paragraphs.Count; $j++) {
    if ($d.Paragraphs.Item($j).Range.Text.Length -eq 1) {
        $blankIndex = $j
        break
    }
}
if ($blankIndex -eq -1) {
    throw "Could not locate the blank paragraph following 'OBS:'"
}

# -----------------------------------------------------------------------
# 1) Remove bold from the paragraph-mark formatting of the "OBS:" paragraph
#    (its pPr/rPr loses <w:b/><w:bCs/>; the runs themselves keep whatever
#    bold formatting they already had).
# -----------------------------------------------------------------------
$obsPara = $d.Paragraphs.Item($obsIndex)
$obsRange = $obsPara.Range
$obsXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t xml:space="preserve">OBS: </w:t></w:r><w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t>M&#233;todos de desenvolvimento onde n&#227;o se &#233; poss&#237;vel retornar, e editar, arquivos que j&#225; foram conclu&#237;dos n&#227;o seriam ideias em vista da arbitrariedade do problema.</w:t></w:r></w:p>
'@
$obsRange.InsertXML($obsXml) | Out-Null

# -----------------------------------------------------------------------
# 2) Replace the first of the four blank paragraphs following "OBS:" with
#    a blank line, a new bold heading question, another blank (bold-styled)
#    line, and a new body paragraph answering the question.
# -----------------------------------------------------------------------
$blankPara = $d.Paragraphs.Item($blankIndex)
$blankRange = $blankPara.Range
$newBlockXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr></w:pPr></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:b/><w:bCs/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:t>O SCRUM deve ser usado no projeto?</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:b/><w:bCs/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr></w:pPr></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:tab/><w:t>Como citado na proposta, n&#227;o h&#225; uma forma fixa de solucionar o problema do sistema de coleta de lixo, por causa das diversas vari&#225;veis envolvidas no ambiente em que a coleta acontece. Portanto, utilizar do SCRUM para o desenvolvimento desde software &#233; extremamente necess&#225;rio. Por o projeto envolver v&#225;rios problemas como tr&#226;nsito, produ&#231;&#227;o do lixo e tempo, &#233; interessante dividir a equipe para cada tarefa realizando reuni&#245;es frequentes para alinhar o que vem sendo produzido por cada uma delas e alterar de planejamento caso algum novo empecilho seja detectado.</w:t></w:r></w:p>
'@
$blankRange.InsertXML($newBlockXml) | Out-Null

Write-Output "done"
